# XLSX export specified worksheet
#
# - Removes the two trailing blank rows from Sheet1.
# - Adds two new sheets ("Another Sheet", "Third Sheet") with their own data.
# - Leaves the Third Sheet as the active / selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1: drop the two empty trailing rows (5 and 6), shrinking the
# used range from A1:E6 down to A1:E4.
$ws1.Rows.Item(5).Resize(2).Delete()

# Add "Another Sheet" right after Sheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Another Sheet"

$ws2.Range("A1").Value = "pink"
$ws2.Range("B1").Value = "green"
$ws2.Range("C1").Value = "blue"
$ws2.Range("A2").Value = "red"
$ws2.Range("B2").Value = "purple"
$ws2.Range("C2").Value = "orange"

# Add "Third Sheet" right after "Another Sheet"; it becomes the active sheet.
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Third Sheet"

$ws3.Range("A1").Value = "something"

$ws3.Range("B6").Select()
